$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style of an untouched, plainly-formatted cell (the sheet default). We
# reapply this after every write below so that forcing a "@" (Text) number
# format -- needed to stop Excel turning numeric-looking strings such as
# "320.82" into real numbers -- does not leave any visible style change
# behind; the cells keep exactly the formatting they started with.
$defaultStyle = $ws.Range("A1").Style

# Cell -> new text value, exactly as described by the commit diff
# (updated crypto prices / 1h volume deltas, plus the Cosmos <-> Toncoin
# row swap at rows 29-30).
$updates = @(
    @{Cell = "D2"; Text = '43.241.43'},
    @{Cell = "E2"; Text = '  -4.95%  '},
    @{Cell = "D3"; Text = '2.240.07'},
    @{Cell = "E3"; Text = '  -5.84%  '},
    @{Cell = "E4"; Text = '  -0.03%  '},
    @{Cell = "D5"; Text = '320.82'},
    @{Cell = "E5"; Text = '  +0.72%  '},
    @{Cell = "D6"; Text = '101.05'},
    @{Cell = "E6"; Text = '  -7.67%  '},
    @{Cell = "E7"; Text = '  -8.43%  '},
    @{Cell = "E8"; Text = '  -0.10%  '},
    @{Cell = "D9"; Text = '0.564'},
    @{Cell = "E9"; Text = '  -8.46%  '},
    @{Cell = "D10"; Text = '37.25'},
    @{Cell = "E10"; Text = '  -9.33%  '},
    @{Cell = "D11"; Text = '54.61'},
    @{Cell = "E11"; Text = '  -2.92%  '},
    @{Cell = "D12"; Text = '0.0830'},
    @{Cell = "E12"; Text = '  -9.87%  '},
    @{Cell = "D13"; Text = '7.72'},
    @{Cell = "E13"; Text = '  -9.90%  '},
    @{Cell = "E14"; Text = '  -0.88%  '},
    @{Cell = "D15"; Text = '2.581.07'},
    @{Cell = "E15"; Text = '  -5.83%  '},
    @{Cell = "D16"; Text = '0.868'},
    @{Cell = "E16"; Text = '  -12.19%  '},
    @{Cell = "E17"; Text = '  -6.27%  '},
    @{Cell = "D18"; Text = '2.244.69'},
    @{Cell = "E18"; Text = '  -5.64%  '},
    @{Cell = "D19"; Text = '43.169.30'},
    @{Cell = "E19"; Text = '  -4.93%  '},
    @{Cell = "D20"; Text = '14.49'},
    @{Cell = "E20"; Text = '  -7.52%  '},
    @{Cell = "E21"; Text = '  -8.80%  '},
    @{Cell = "D22"; Text = '6.54'},
    @{Cell = "E22"; Text = '  -10.88%  '},
    @{Cell = "D23"; Text = '65.66'},
    @{Cell = "E23"; Text = '  -10.52%  '},
    @{Cell = "E24"; Text = '  -14.16%  '},
    @{Cell = "E25"; Text = '  -8.78%  '},
    @{Cell = "E26"; Text = '  -7.12%  '},
    @{Cell = "D27"; Text = '1.00'},
    @{Cell = "E27"; Text = '  -0.20%  '},
    @{Cell = "D28"; Text = '4.04'},
    @{Cell = "E28"; Text = '  +0.94%  '},
    @{Cell = "B29"; Text = 'Toncoin'},
    @{Cell = "C29"; Text = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'},
    @{Cell = "D29"; Text = '2.25'},
    @{Cell = "E29"; Text = '  -1.84%  '},
    @{Cell = "B30"; Text = 'Cosmos'},
    @{Cell = "C30"; Text = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Cell = "D30"; Text = '10.06'},
    @{Cell = "E30"; Text = '  -10.92%  '},
    @{Cell = "D31"; Text = '6.38'},
    @{Cell = "E31"; Text = '  -16.13%  '},
    @{Cell = "D32"; Text = '36.41'},
    @{Cell = "E32"; Text = '  -3.28%  '},
    @{Cell = "D33"; Text = '0.0882'},
    @{Cell = "E33"; Text = '  -7.43%  '},
    @{Cell = "E34"; Text = '  -8.99%  '},
    @{Cell = "D35"; Text = '152.32'},
    @{Cell = "E35"; Text = '  -9.59%  '},
    @{Cell = "E36"; Text = '  -6.41%  '},
    @{Cell = "D37"; Text = '3.23'},
    @{Cell = "E37"; Text = '  +6.48%  '},
    @{Cell = "D38"; Text = '1.96'},
    @{Cell = "E38"; Text = '  +0.46%  '},
    @{Cell = "E39"; Text = '  -7.96%  '},
    @{Cell = "E40"; Text = '  -6.38%  '},
    @{Cell = "D41"; Text = '0.104'},
    @{Cell = "E41"; Text = '  -10.70%  '},
    @{Cell = "D42"; Text = '3.67'},
    @{Cell = "E42"; Text = '  -8.42%  '},
    @{Cell = "D43"; Text = '0.0325'},
    @{Cell = "E43"; Text = '  -8.80%  '},
    @{Cell = "D44"; Text = '13.86'},
    @{Cell = "E44"; Text = '  +6.69%  '},
    @{Cell = "E45"; Text = '  -0.02%  '},
    @{Cell = "D46"; Text = '1.751.76'},
    @{Cell = "E46"; Text = '  -5.42%  '},
    @{Cell = "D47"; Text = '86.69'},
    @{Cell = "E48"; Text = '  -10.29%  '},
    @{Cell = "D49"; Text = '5.35'},
    @{Cell = "E49"; Text = '  -10.83%  '},
    @{Cell = "D50"; Text = '76.00'},
    @{Cell = "E50"; Text = '  -9.36%  '},
    @{Cell = "D51"; Text = '59.10'},
    @{Cell = "E51"; Text = '  -16.25%  '}
)

foreach ($update in $updates) {
    $rng = $ws.Range($update.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $update.Text
    $rng.Style = $defaultStyle
}
